$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C21").NumberFormat = "@"
$ws.Range("C2").Value = "111110000001011000010101000100"
$ws.Range("D2").Value = 0.9391295160492703
$ws.Range("E2").Value = 0.00487892865040852
$ws.Range("F2").Value = 0.2967567476953556
$ws.Range("C3").Value = "111110000100100000001111101101"
$ws.Range("D3").Value = 0.9406081878374415
$ws.Range("E3").Value = 0.3014185174356079
$ws.Range("F3").Value = 0.5828704773098387
$ws.Range("C4").Value = "111110000100100000001111101101"
$ws.Range("D4").Value = 0.9406081878374415
$ws.Range("E4").Value = 0.5999076720498344
$ws.Range("F4").Value = 0.9059385340677079
$ws.Range("C5").Value = "111110000100100000001111101101"
$ws.Range("D5").Value = 0.9406081878374415
$ws.Range("E5").Value = 0.9391307452990107
$ws.Range("F5").Value = 0.9404522517087889
$ws.Range("C6").Value = "111110000100100000001111101101"
$ws.Range("D6").Value = 0.9406081878374415
$ws.Range("E6").Value = 0.9406081137715335
$ws.Range("F6").Value = 0.9406081804308506
$ws.Range("C7").Value = "111110000100100000001111101101"
$ws.Range("D7").Value = 0.9406081878374415
$ws.Range("E7").Value = 0.9406081878374415
$ws.Range("F7").Value = 0.9406081878374414
$ws.Range("C8").Value = "111110000100100000001111101101"
$ws.Range("D8").Value = 0.9406081878374415
$ws.Range("E8").Value = 0.9406081878374415
$ws.Range("F8").Value = 0.9406081878374414
$ws.Range("C9").Value = "111110000100100000001111101101"
$ws.Range("D9").Value = 0.9406081878374415
$ws.Range("E9").Value = 0.9406081878374415
$ws.Range("F9").Value = 0.9406081878374414
$ws.Range("C10").Value = "111110000100100000001111101101"
$ws.Range("D10").Value = 0.9406081878374415
$ws.Range("E10").Value = 0.9406081878374415
$ws.Range("F10").Value = 0.9406081878374414
$ws.Range("C11").Value = "111110000100100000001111101101"
$ws.Range("D11").Value = 0.9406081878374415
$ws.Range("E11").Value = 0.9406081878374415
$ws.Range("F11").Value = 0.9406081878374414
$ws.Range("C12").Value = "111110000100100000001111101101"
$ws.Range("D12").Value = 0.9406081878374415
$ws.Range("E12").Value = 0.9406081878374415
$ws.Range("F12").Value = 0.9406081878374414
$ws.Range("C13").Value = "111110000110100000001111101101"
$ws.Range("D13").Value = 0.9415555449773083
$ws.Range("E13").Value = 0.9406077253771852
$ws.Range("F13").Value = 0.9407028773054025
$ws.Range("C14").Value = "111110000110100000011111101101"
$ws.Range("D14").Value = 0.9415573957507929
$ws.Range("E14").Value = 0.2207586177049297
$ws.Range("F14").Value = 0.8634143462895574
$ws.Range("C15").Value = "111110000110100000011111101101"
$ws.Range("D15").Value = 0.9415573957507929
$ws.Range("E15").Value = 0.9415555449773083
$ws.Range("F15").Value = 0.9415562852867023
$ws.Range("C16").Value = "111110000110100000011111101101"
$ws.Range("D16").Value = 0.9415573957507929
$ws.Range("E16").Value = 0.9415573957507929
$ws.Range("F16").Value = 0.9415573957507929
$ws.Range("C17").Value = "111110000110100000011111101101"
$ws.Range("D17").Value = 0.9415573957507929
$ws.Range("E17").Value = 0.9415572800773967
$ws.Range("F17").Value = 0.9415573841834532
$ws.Range("C18").Value = "111110000110100000011111101101"
$ws.Range("D18").Value = 0.9415573957507929
$ws.Range("E18").Value = 0.9415573957507929
$ws.Range("F18").Value = 0.9415573957507929
$ws.Range("C19").Value = "111110000110100000011111101101"
$ws.Range("D19").Value = 0.9415573957507929
$ws.Range("E19").Value = 0.9415573957507929
$ws.Range("F19").Value = 0.9415573957507929
$ws.Range("C20").Value = "111110000110100000011111101101"
$ws.Range("D20").Value = 0.9415573957507929
$ws.Range("E20").Value = 0.9415573957507929
$ws.Range("F20").Value = 0.9415573957507929
$ws.Range("C21").Value = "111110000110100000011111101101"
$ws.Range("D21").Value = 0.9415573957507929
$ws.Range("E21").Value = 0.941557381291618
$ws.Range("F21").Value = 0.9415573943048756
